{"js": "// Add a new \"Compact List\" paragraph style (styleId \"CompactList\") that\n// mirrors the existing \"Compact\" style: based on \"Body Text\", a Quick\n// Style, with 36-twip (1.8pt) space before/after.\ncontext.document.addStyle(\"Compact List\", Word.StyleType.paragraph);\nawait context.sync();\n\n// addStyle()'s own return value isn't reliably bound to the freshly minted\n// style, so re-fetch it from the style collection before configuring it.\nconst style = context.document.getStyles().getByName(\"Compact List\");\nstyle.baseStyle = \"BodyText\";\nstyle.quickStyle = true;\nstyle.paragraphFormat.spaceBefore = 1.8;\nstyle.paragraphFormat.spaceAfter = 1.8;\nawait context.sync();\n", "ps1": "# Add a new \"Compact List\" paragraph style (styleId \"CompactList\") that\n# mirrors the existing \"Compact\" style: based on \"Body Text\", a Quick\n# Style, with 36-twip (1.8pt) space before/after.\n$d = $word.ActiveDocument\n\n# 1 = wdStyleTypeParagraph\n$style = $d.Styles.Add(\"Compact List\", 1)\n$style.BaseStyle = \"BodyText\"\n$style.QuickStyle = $true\n$style.ParagraphFormat.SpaceBefore = 1.8\n$style.ParagraphFormat.SpaceAfter = 1.8\n"}
